# Apply the artfynd "A 50292-2021" update:
#  - round the Ost/Nord coordinates on row 13 to whole metres
#  - drop the (empty) Starttid/Sluttid cells on row 13
#  - append two new observation rows (14 and 15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 edits -----------------------------------------------------
$ws.Cells.Item(13, 17).Value = 478586          # Q13 Ost
$ws.Cells.Item(13, 18).Value = 6556137         # R13 Nord

$ws.Cells.Item(13, 26).ClearContents()         # Z13 Starttid -> removed
$ws.Cells.Item(13, 28).ClearContents()         # AB13 Sluttid -> removed

# --- Row 14: Stubbspretmossa (Herzogiella seligeri) --------------------
$ws.Cells.Item(14, 1).Value = 112231588
$ws.Cells.Item(14, 2).Value = 93171
$ws.Cells.Item(14, 3).Value = "Ovaliderad"
$ws.Cells.Item(14, 4).Value = "LC"
$ws.Cells.Item(14, 5).Value = 2818
$ws.Cells.Item(14, 6).Value = "Stubbspretmossa"
$ws.Cells.Item(14, 7).Value = "Herzogiella seligeri"
$ws.Cells.Item(14, 8).Value = "(Brid.) Z.Iwats."

$ws.Cells.Item(14, 16).Value = "Baggetorp, Nrk"
$ws.Cells.Item(14, 17).Value = 478719
$ws.Cells.Item(14, 18).Value = 6556487
$ws.Cells.Item(14, 19).Value = 10
$ws.Cells.Item(14, 20).Value = "Örebro"
$ws.Cells.Item(14, 21).Value = "Lekeberg"
$ws.Cells.Item(14, 22).Value = "Närke"
$ws.Cells.Item(14, 23).Value = "Kvistbro"

$ws.Cells.Item(14, 25).NumberFormat = "@"
$ws.Cells.Item(14, 25).Value = "2023-09-19"    # Y14 Startdatum
$ws.Cells.Item(14, 27).NumberFormat = "@"
$ws.Cells.Item(14, 27).Value = "2023-09-19"    # AA14 Slutdatum

$ws.Cells.Item(14, 30).Value = $false
$ws.Cells.Item(14, 31).Value = $false
$ws.Cells.Item(14, 33).Value = $false
$ws.Cells.Item(14, 34).Value = "Blåbärsbarrskog"

$ws.Cells.Item(14, 36).Value = "tall"
$ws.Cells.Item(14, 37).Value = "Pinus sylvestris"
$ws.Cells.Item(14, 39).Value = "Liggande död trädstam, markontakt"
$ws.Cells.Item(14, 41).Value = "Horizontal, dead with ground contact # murken grov låga # Pinus sylvestris"

$ws.Cells.Item(14, 49).Value = "Michael Andersson"
$ws.Cells.Item(14, 50).Value = "Michael Andersson"

# --- Row 15: Talltita (Poecile montanus) --------------------------------
$ws.Cells.Item(15, 1).Value = 112231491
$ws.Cells.Item(15, 2).Value = 56543
$ws.Cells.Item(15, 3).Value = "Ovaliderad"
$ws.Cells.Item(15, 4).Value = "NT"
$ws.Cells.Item(15, 5).Value = 103021
$ws.Cells.Item(15, 6).Value = "Talltita"
$ws.Cells.Item(15, 7).Value = "Poecile montanus"
$ws.Cells.Item(15, 8).Value = "(Conrad von Baldenstein, 1827)"

$ws.Cells.Item(15, 9).NumberFormat = "@"
$ws.Cells.Item(15, 9).Value = "1"              # I15 Antal (kept as text)

$ws.Cells.Item(15, 13).Value = "lockläte, övriga läten"

$ws.Cells.Item(15, 16).Value = "Baggetorp, Nrk"
$ws.Cells.Item(15, 17).Value = 478579
$ws.Cells.Item(15, 18).Value = 6556322
$ws.Cells.Item(15, 19).Value = 10
$ws.Cells.Item(15, 20).Value = "Örebro"
$ws.Cells.Item(15, 21).Value = "Lekeberg"
$ws.Cells.Item(15, 22).Value = "Närke"
$ws.Cells.Item(15, 23).Value = "Kvistbro"

$ws.Cells.Item(15, 25).NumberFormat = "@"
$ws.Cells.Item(15, 25).Value = "2023-09-19"    # Y15 Startdatum
$ws.Cells.Item(15, 27).NumberFormat = "@"
$ws.Cells.Item(15, 27).Value = "2023-09-19"    # AA15 Slutdatum

$ws.Cells.Item(15, 30).Value = $false
$ws.Cells.Item(15, 31).Value = $false
$ws.Cells.Item(15, 33).Value = $false

$ws.Cells.Item(15, 49).Value = "Michael Andersson"
$ws.Cells.Item(15, 50).Value = "Michael Andersson"
